# 2017-01-31 update: energy.gov - chunk 7
# Table 4.12.A - roll the report forward from October 2016/2015 to November 2016/2015
# and refresh the associated state-level data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title & column header text: October -> November -------------------
$ws.Range("A1").Value = "Table 4.12.A. Average Cost of Petroleum Coke Delivered for Electricity Generation by State, November 2016 and 2015"

# The header cells (B4/C4/E4/F4/G4/H4) carry a numeric cell style even
# though they hold text labels like "November 2016" - Excel's
# autorecognition would otherwise read that text as a date. Force the
# cell to Text format while the value is written, then restore the
# original numeric format so the cell style/appearance is unchanged.
$headerCells = "B4", "C4", "E4", "F4", "G4", "H4"
$headerValues = "November 2016", "November 2015", "November 2016", "November 2015", "November 2016", "November 2015"
for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $cell = $ws.Range($headerCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $headerValues[$i]
    $cell.NumberFormat = "0.00"
}

# --- Data updates ---------------------------------------------------------

# Row 16: East North Central
$ws.Range("E16").Value = 1.45
$ws.Range("F16").Value = 1.1399999999999999

# Row 19: Michigan
$ws.Range("B19").Value = 1.32
$ws.Range("D19").Value = -0.2
$ws.Range("E19").Value = 1.32

# Row 21: Wisconsin
$ws.Range("B21").Value = 1.73
$ws.Range("C21").Value = 1.73
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1.73
$ws.Range("F21").Value = 1.73

# Row 30: South Atlantic
$ws.Range("B30").Value = 2.35
$ws.Range("C30").Value = 1.65
$ws.Range("D30").Value = 0.42
$ws.Range("E30").Value = 2.35
$ws.Range("F30").Value = 1.65

# Row 33: Florida
$ws.Range("B33").Value = 2.35
$ws.Range("C33").Value = 1.65
$ws.Range("D33").Value = 0.42
$ws.Range("E33").Value = 2.35
$ws.Range("F33").Value = 1.65

# Row 40: East South Central
$ws.Range("B40").Value = 1.47
$ws.Range("C40").Value = 1.57
$ws.Range("D40").Value = -0.064
$ws.Range("E40").Value = 1.47
$ws.Range("F40").Value = 1.57

# Row 42: Kentucky
$ws.Range("B42").Value = 1.47
$ws.Range("C42").Value = 1.57
$ws.Range("D42").Value = -0.064
$ws.Range("E42").Value = 1.47
$ws.Range("F42").Value = 1.57

# Row 45: West South Central
$ws.Range("B45").Value = 2.34
$ws.Range("C45").Value = 1.49
$ws.Range("D45").Value = 0.57
$ws.Range("E45").Value = 2.34
$ws.Range("F45").Value = 1.49

# Row 47: Louisiana
$ws.Range("B47").Value = 2.34
$ws.Range("C47").Value = 1.49
$ws.Range("D47").Value = 0.57
$ws.Range("E47").Value = 2.34
$ws.Range("F47").Value = 1.49

# Row 66: U.S. Total
$ws.Range("E66").Value = 2.2200000000000002
$ws.Range("F66").Value = 1.46
